$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure target cells remain text (avoid Excel auto-converting numeric-looking
# or percent-looking strings into numbers/percentages).
$targetRows = @(2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,38,39,40,41,42,43,44,45,46,47,48,49,50)
foreach ($r in $targetRows) {
    $ws.Range("D$r").NumberFormat = "@"
    $ws.Range("E$r").NumberFormat = "@"
}

# Apply updated price (D) and volume/1h change (E) values
$ws.Range("D2").Value = "327.14"
$ws.Range("E2").Value = "0.21%"
$ws.Range("D3").Value = "43.94"
$ws.Range("E3").Value = "-2.13%"
$ws.Range("D4").Value = "5.522"
$ws.Range("E4").Value = "-0.59%"
$ws.Range("D5").Value = "0.08025"
$ws.Range("E5").Value = "-0.66%"
$ws.Range("D6").Value = "1.990"
$ws.Range("E6").Value = "4.63%"
$ws.Range("D7").Value = "2.574"
$ws.Range("E7").Value = "-4.92%"
$ws.Range("D8").Value = "0.9530"
$ws.Range("E8").Value = "0.51%"
$ws.Range("D9").Value = "0.1130"
$ws.Range("E9").Value = "-4.04%"
$ws.Range("D10").Value = "0.1864"
$ws.Range("E10").Value = "-1.50%"
$ws.Range("D11").Value = "10.66"
$ws.Range("E11").Value = "23.14%"
$ws.Range("D12").Value = "0.09842"
$ws.Range("E12").Value = "-2.84%"
$ws.Range("D13").Value = "0.04598"
$ws.Range("E13").Value = "9.69%"
$ws.Range("D14").Value = "0.1066"
$ws.Range("E14").Value = "0.11%"
$ws.Range("D15").Value = "0.001262"
$ws.Range("E15").Value = "-1.51%"
$ws.Range("D16").Value = "0.04090"
$ws.Range("E16").Value = "-3.73%"
$ws.Range("D17").Value = "0.005862"
$ws.Range("E17").Value = "-2.01%"
$ws.Range("E18").Value = "-6.66%"
$ws.Range("D19").Value = "4.296"
$ws.Range("E19").Value = "-0.91%"
$ws.Range("D20").Value = "0.3477"
$ws.Range("E20").Value = "-0.26%"
$ws.Range("D21").Value = "0.1407"
$ws.Range("E21").Value = "2.47%"
$ws.Range("E22").Value = "-4.42%"
$ws.Range("E23").Value = "0.65%"
$ws.Range("D24").Value = "0.004335"
$ws.Range("E24").Value = "-6.43%"
$ws.Range("D25").Value = "0.0001189"
$ws.Range("E25").Value = "-3.51%"
$ws.Range("D26").Value = "0.0003740"
$ws.Range("E26").Value = "-6.51%"
$ws.Range("D38").Value = "0.02550"
$ws.Range("E38").Value = "-4.05%"
$ws.Range("D39").Value = "0.05650"
$ws.Range("E39").Value = "1.73%"
$ws.Range("D40").Value = "0.007534"
$ws.Range("E40").Value = "-1.97%"
$ws.Range("D41").Value = "0.1396"
$ws.Range("E41").Value = "0.16%"
$ws.Range("D42").Value = "0.007588"
$ws.Range("E42").Value = "-33.06%"
$ws.Range("D43").Value = "0.002013"
$ws.Range("E43").Value = "-2.18%"
$ws.Range("D44").Value = "0.008493"
$ws.Range("E44").Value = "-7.88%"
$ws.Range("D45").Value = "0.00007103"
$ws.Range("E45").Value = "-0.15%"
$ws.Range("E46").Value = "-0.40%"
$ws.Range("D47").Value = "0.003529"
$ws.Range("E47").Value = "55.04%"
$ws.Range("D48").Value = "0.003028"
$ws.Range("E48").Value = "-12.06%"
$ws.Range("D49").Value = "0.00002097"
$ws.Range("E49").Value = "-0.40%"
$ws.Range("D50").Value = "0.0001998"
$ws.Range("E50").Value = "-0.40%"
